# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2404"
#   "<header>_new" -> "<header>_FV2410"
# Wrap the data range in an Excel Table (ListObject) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header base names, in the order they appear left to right.
$headerNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J (1-10) carry the "_old" -> "_FV2404" headers.
for ($i = 0; $i -lt $headerNames.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($headerNames[$i])_FV2404"
}

# Column K (11) is the "diff" column and stays untouched.

# Columns L..U (12-21) carry the "_new" -> "_FV2410" headers.
for ($i = 0; $i -lt $headerNames.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($headerNames[$i])_FV2410"
}

# Turn the whole data range into an Excel table ("Table1") with the header row.
$dataRange = $ws.Range("A1:U79")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# Freeze the header row (split after row 1, keep the header pane visible while scrolling).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
